# chore: update Sheets via scheduled runner
# Refreshes currentAveragePrice*/LevePrice*/LeveProfit* market-board figures
# across the ALC/ARM/BSM/CRP/CUL/GSM/LTW/WVR leve-profit tables.
$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H40").Value = 2902.647
$ws.Range("I40").Value = 1920
$ws.Range("J40").Value = 3312.0833
$ws.Range("K40").Value = 1920
$ws.Range("L40").Value = 3312.0833
$ws.Range("M40").Value = -1745
$ws.Range("N40").Value = -3662.0833
$ws.Range("H58").Value = 103.333336
$ws.Range("J58").Value = 0
$ws.Range("L58").Value = 0
$ws.Range("H112").Value = 1364.2
$ws.Range("J112").Value = 1388.6316
$ws.Range("L112").Value = 4165.8948
$ws.Range("N112").Value = -6381.8948
$ws.Range("H129").Value = 1527.1875
$ws.Range("I129").Value = 1404.3636
$ws.Range("K129").Value = 4213.0908
$ws.Range("M129").Value = 786.9092000000001
$ws.Range("H132").Value = 11245.857
$ws.Range("I132").Value = 12894.111
$ws.Range("J132").Value = 8279
$ws.Range("K132").Value = 38682.333
$ws.Range("L132").Value = 24837
$ws.Range("M132").Value = -36152.333
$ws.Range("N132").Value = -29897
$ws.Range("H138").Value = 1672.9574
$ws.Range("I138").Value = 1216.4348
$ws.Range("J138").Value = 2110.4583
$ws.Range("K138").Value = 3649.3044
$ws.Range("L138").Value = 6331.374899999999
$ws.Range("M138").Value = 1490.6956
$ws.Range("N138").Value = -16611.3749
$ws.Range("N58").ClearContents()

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H32").Value = 7608.4473
$ws.Range("I32").Value = 3942.8708
$ws.Range("J32").Value = 23841.715
$ws.Range("K32").Value = 3942.8708
$ws.Range("L32").Value = 23841.715
$ws.Range("M32").Value = -3655.8708
$ws.Range("N32").Value = -24415.715
$ws.Range("H88").Value = 8773472
$ws.Range("I88").Value = 23810424
$ws.Range("J88").Value = 1916.6666
$ws.Range("K88").Value = 23810424
$ws.Range("L88").Value = 1916.6666
$ws.Range("M88").Value = -23810018
$ws.Range("N88").Value = -2728.6666
$ws.Range("H91").Value = 8773472
$ws.Range("I91").Value = 23810424
$ws.Range("J91").Value = 1916.6666
$ws.Range("K91").Value = 23810424
$ws.Range("L91").Value = 1916.6666
$ws.Range("M91").Value = -23809020
$ws.Range("N91").Value = -4724.6666
$ws.Range("H139").Value = 59499.5
$ws.Range("J139").Value = 59499.5
$ws.Range("L139").Value = 59499.5
$ws.Range("N139").Value = -69779.5

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H81").Value = 0
$ws.Range("J81").Value = 0
$ws.Range("L81").Value = 0
$ws.Range("H84").Value = 0
$ws.Range("J84").Value = 0
$ws.Range("L84").Value = 0
$ws.Range("H86").Value = 2362.4
$ws.Range("I86").Value = 1680
$ws.Range("J86").Value = 3196.4443
$ws.Range("K86").Value = 1680
$ws.Range("L86").Value = 3196.4443
$ws.Range("M86").Value = -557
$ws.Range("N86").Value = -5442.4443
$ws.Range("H89").Value = 2362.4
$ws.Range("I89").Value = 1680
$ws.Range("J89").Value = 3196.4443
$ws.Range("K89").Value = 8400
$ws.Range("L89").Value = 15982.2215
$ws.Range("M89").Value = -2784
$ws.Range("N89").Value = -27214.2215
$ws.Range("H105").Value = 1853.5405
$ws.Range("J105").Value = 2261
$ws.Range("L105").Value = 2261
$ws.Range("N105").Value = -5755
$ws.Range("N81").ClearContents()
$ws.Range("N84").ClearContents()

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H10").Value = 565
$ws.Range("I10").Value = 565
$ws.Range("J10").Value = 0
$ws.Range("K10").Value = 565
$ws.Range("L10").Value = 0
$ws.Range("M10").Value = -426
$ws.Range("H51").Value = 55000
$ws.Range("J51").Value = 55000
$ws.Range("L51").Value = 55000
$ws.Range("N51").Value = -56472
$ws.Range("H61").Value = 55000
$ws.Range("J61").Value = 55000
$ws.Range("L61").Value = 55000
$ws.Range("N61").Value = -55696
$ws.Range("H62").Value = 129251.25
$ws.Range("I62").Value = 3502.5
$ws.Range("K62").Value = 3502.5
$ws.Range("M62").Value = -2878.5
$ws.Range("H65").Value = 129251.25
$ws.Range("I65").Value = 3502.5
$ws.Range("K65").Value = 17512.5
$ws.Range("M65").Value = -14392.5
$ws.Range("H88").Value = 11187.223
$ws.Range("J88").Value = 11187.223
$ws.Range("L88").Value = 11187.223
$ws.Range("N88").Value = -11999.223
$ws.Range("H91").Value = 11187.223
$ws.Range("J91").Value = 11187.223
$ws.Range("L91").Value = 11187.223
$ws.Range("N91").Value = -13995.223
$ws.Range("H105").Value = 1068.4445
$ws.Range("I105").Value = 1262.5
$ws.Range("J105").Value = 913.2
$ws.Range("K105").Value = 1262.5
$ws.Range("L105").Value = 913.2
$ws.Range("M105").Value = 484.5
$ws.Range("N105").Value = -4407.2
$ws.Range("H122").Value = 1179.4348
$ws.Range("J122").Value = 1833.3334
$ws.Range("L122").Value = 5500.0002
$ws.Range("N122").Value = -10400.0002
$ws.Range("H132").Value = 3003.5
$ws.Range("I132").Value = 3003.5
$ws.Range("K132").Value = 9010.5
$ws.Range("M132").Value = -6480.5
$ws.Range("H134").Value = 6759.375
$ws.Range("I134").Value = 4753
$ws.Range("K134").Value = 14259
$ws.Range("M134").Value = -11724
$ws.Range("N10").ClearContents()

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H4").Value = 16367651
$ws.Range("I4").Value = 16522959
$ws.Range("K4").Value = 49568877
$ws.Range("M4").Value = -49568765
$ws.Range("H36").Value = 8666.666999999999
$ws.Range("I36").Value = 500
$ws.Range("K36").Value = 1500
$ws.Range("M36").Value = -1331
$ws.Range("H137").Value = 2043.1578
$ws.Range("I137").Value = 1178.625
$ws.Range("J137").Value = 2671.9092
$ws.Range("K137").Value = 3535.875
$ws.Range("L137").Value = 8015.7276
$ws.Range("M137").Value = 1564.125
$ws.Range("N137").Value = -18215.7276

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H92").Value = 8199.166999999999
$ws.Range("J92").Value = 8199.166999999999
$ws.Range("L92").Value = 8199.166999999999
$ws.Range("N92").Value = -11943.167
$ws.Range("H102").Value = 1492.9623
$ws.Range("I102").Value = 706.7368
$ws.Range("K102").Value = 706.7368
$ws.Range("M102").Value = 915.2632
$ws.Range("H113").Value = 2925.4666
$ws.Range("I113").Value = 1147.1666
$ws.Range("J113").Value = 4111
$ws.Range("K113").Value = 1147.1666
$ws.Range("L113").Value = 4111
$ws.Range("M113").Value = 1022.8334
$ws.Range("N113").Value = -8451
$ws.Range("H122").Value = 2588.2778
$ws.Range("I122").Value = 2208.5454
$ws.Range("J122").Value = 3185
$ws.Range("K122").Value = 6625.6362
$ws.Range("L122").Value = 9555
$ws.Range("M122").Value = -4175.6362
$ws.Range("N122").Value = -14455

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H40").Value = 2965.6365
$ws.Range("I40").Value = 2928.2856
$ws.Range("K40").Value = 2928.2856
$ws.Range("M40").Value = -2792.2856
$ws.Range("H46").Value = 2602.1052
$ws.Range("I46").Value = 1800
$ws.Range("J46").Value = 2696.4707
$ws.Range("K46").Value = 1800
$ws.Range("L46").Value = 2696.4707
$ws.Range("M46").Value = -1612
$ws.Range("N46").Value = -3072.4707
$ws.Range("H55").Value = 601.43475
$ws.Range("J55").Value = 561.46155
$ws.Range("L55").Value = 561.46155
$ws.Range("N55").Value = -907.46155
$ws.Range("H61").Value = 11801.333
$ws.Range("I61").Value = 13161.6
$ws.Range("K61").Value = 13161.6
$ws.Range("M61").Value = -12959.6
$ws.Range("H93").Value = 10103467
$ws.Range("I93").Value = 12822877
$ws.Range("K93").Value = 12822877
$ws.Range("M93").Value = -12821629
$ws.Range("H113").Value = 11801.333
$ws.Range("I113").Value = 13161.6
$ws.Range("K113").Value = 13161.6
$ws.Range("M113").Value = -10991.6
$ws.Range("H119").Value = 45000
$ws.Range("J119").Value = 45000
$ws.Range("L119").Value = 45000
$ws.Range("N119").Value = -54676
$ws.Range("H122").Value = 4071.4285
$ws.Range("I122").Value = 4500
$ws.Range("K122").Value = 13500
$ws.Range("M122").Value = -11050

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H18").Value = 10999
$ws.Range("I18").Value = 10999
$ws.Range("K18").Value = 10999
$ws.Range("M18").Value = -10826
$ws.Range("H43").Value = 37125
$ws.Range("I43").Value = 38428.57
$ws.Range("K43").Value = 38428.57
$ws.Range("M43").Value = -38279.57
$ws.Range("H81").Value = 12503574
$ws.Range("I81").Value = 2419.1
$ws.Range("K81").Value = 4838.2
$ws.Range("M81").Value = -3777.2
$ws.Range("H84").Value = 12503574
$ws.Range("I84").Value = 2419.1
$ws.Range("K84").Value = 24191
$ws.Range("M84").Value = -18887
$ws.Range("H92").Value = 80000
$ws.Range("J92").Value = 80000
$ws.Range("L92").Value = 80000
$ws.Range("N92").Value = -84992
$ws.Range("H107").Value = 771.4211
$ws.Range("I107").Value = 656.8125
$ws.Range("K107").Value = 1970.4375
$ws.Range("M107").Value = -50.4375
